# Applies the "yolo & rcnn assignment complete" edit:
#   1. Assignment Name -> "Assignment Name -  CNN Architecture" (with
#      proofErr gramStart/gramEnd markers around "-  CNN", as Word's
#      grammar checker would insert for a sentence fragment starting
#      with a dash).
#   2. Submission Date -> day changes from 11 to 20 (11-12-2024 -> 20-12-2024).
#   3. Git Link -> points at the new repo path / notebook file name.
#
# Each paragraph is rewritten in place via Range.InsertXML so the
# surrounding <w:p> (paraId / rsid / pPr) is preserved and only the
# run content inside changes, mirroring the target OOXML diff.

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$d = $word.ActiveDocument

function Set-ParagraphRuns($paragraph, [string]$innerXml) {
    $start = $paragraph.Range.Start
    $end = $paragraph.Range.End - 1   # exclude the trailing paragraph mark
    $range = $d.Range($start, $end)
    $range.InsertXML("<w:p $wordNs>" + $innerXml + "</w:p>")
}

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -like "Assignment Name*") {
        $inner = '<w:r><w:t xml:space="preserve">Assignment Name </w:t></w:r>' +
                 '<w:proofErr w:type="gramStart"/>' +
                 '<w:r><w:t xml:space="preserve">-  </w:t></w:r>' +
                 '<w:r><w:t>CNN</w:t></w:r>' +
                 '<w:proofErr w:type="gramEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                 '<w:r><w:t>Architecture</w:t></w:r>'
        Set-ParagraphRuns $p $inner
    }
    elseif ($t -like "Submission Date*") {
        $inner = '<w:r><w:t xml:space="preserve">Submission Date &#8211; </w:t></w:r>' +
                 '<w:r><w:t>20</w:t></w:r>' +
                 '<w:r><w:t>-</w:t></w:r>' +
                 '<w:r><w:t>12</w:t></w:r>' +
                 '<w:r><w:t xml:space="preserve">-2024 </w:t></w:r>'
        Set-ParagraphRuns $p $inner
    }
    elseif ($t -like "Git Link*") {
        $inner = '<w:r><w:t xml:space="preserve">Git Link </w:t></w:r>' +
                 '<w:r><w:t>-</w:t></w:r>' +
                 '<w:r><w:t>https://github.com/ankitsharma5911/deeplearning-assignment/blob/main/CNN%20Architecture.ipynb</w:t></w:r>'
        Set-ParagraphRuns $p $inner
    }
}
